# Issue #78 - Added descriptions to cards
# Populate the "Description" column (F) on the Activities sheet with the
# card description text, matching each card's Title, and fix the database
# setup script: the "Multiple perspectives" event's save condition on the
# Events sheet was updated from ">1" to ">=2".

$wb = $excel.ActiveWorkbook

$act = $wb.Worksheets.Item("Activities")

$act.Range("F2").Value  = "Contingency time is time allocated during planning for unforeseen events or problems that may arise."
$act.Range("F3").Value  = "Do you know what the ethical clearance procedure is for your project? If not you need to find out – quickly!"
$act.Range("F4").Value  = "This refers to recognizing where you will need additional training and taking steps to learn what you need to before your project starts."
$act.Range("F5").Value  = "Holidays and breaks are really important. But don’t take too many!"
$act.Range("F6").Value  = "Milestones are internal deadlines that you plan to keep yourself on track throughout the project, for example, completing a chapter, releasing a survey, or finishing a product prototype."
$act.Range("F7").Value  = "This means looking after your data and metadata. As a minimum this could simply be ensuring that you backup your work regularly, or it could be much more complex digital curation techniques."
$act.Range("F8").Value  = "Your supervisor is a precious resource, plan your time with them well!"
$act.Range("F9").Value  = "Thinking is great! But it doesn’t directly produce anything. (Watch out for no upwards arrow!) Too much thinking and not enough planning will hamper your progress."
$act.Range("F10").Value = "Other people are a great contextual resource!"
$act.Range("F11").Value = "Sometimes an article that’s not relevant to you can still lead you to other good sources of information. (Watch out for no upwards arrow!)"
$act.Range("F12").Value = "A methodology is a system of methods of doing or studying something. It is part of your design process for carrying out research and investigating others’ methods and methodologies is part of contextual review."
$act.Range("F13").Value = "A methodology is a system of methods of doing or studying something. It is part of your design process for carrying out research and investigating others’ methods and methodologies is part of contextual review."
$act.Range("F14").Value = "You may want to tweak your research questions as you do your contextual research, this will help you focus your study."

# F15 also gets a small font-format tweak (11pt Calibri, black) on top of the text.
$act.Range("F15").Value = "A reference management system will quickly pay back the time you take to learn how to use it. If you don’t already know of one, ask your tutor."
$f15Font = $act.Range("F15").Font
$f15Font.Name = "Calibri"
$f15Font.Size = 11
$f15Font.Color = 0

$act.Range("F20").Value = "This could be an existing database you have accessed and will use directly in your research."
$act.Range("F21").Value = "This could be an existing database you have accessed and will use directly in your research."
$act.Range("F26").Value = "These could be physical or digital objects that you are studying or using in your research, for example artworks or existing apps."
$act.Range("F27").Value = "These could be physical or digital objects that you are studying or using in your research, for example artworks or existing apps."
$act.Range("F28").Value = "This could be a government or industrial policy that is relevant to your research."
$act.Range("F29").Value = "This could be a government or industrial policy that is relevant to your research."
$act.Range("F33").Value = "A methodology is a system of methods of doing or studying something. It is part of your design process for carrying out research and investigating others’ methods and methodologies is part of contextual review."
$act.Range("F51").Value = "Descriptive writing is necessary but it needs to be followed up by some explanation and critical analysis so that readers understand your research in context."
$act.Range("F52").Value = "Descriptive writing is necessary but it needs to be followed up by some explanation and critical analysis so that readers understand your research in context."
$act.Range("F53").Value = "Descriptive writing is necessary but it needs to be followed up by some explanation and critical analysis so that readers understand your research in context."
$act.Range("F54").Value = "Evaluating your work helps you to evidence your conclusions. Make sure you critically analyse your findings rather than just describing them."
$act.Range("F57").Value = "When writing up your conclusions, revisit your research questions and check that you are actually answering them!"
$act.Range("F59").Value = "Although not part of what you submit, having support from family and/or friends can be really important for your wellbeing as you progress your research!"

# Fixed the database setup script: "Multiple perspectives" event save
# condition tightened from ">1" to ">=2" (needs 2+ Resource cards, not just 1).
$events = $wb.Worksheets.Item("Events")
$events.Range("F15").Value = "(17+18+19+20+21+22+23+24+25+26+27+28):>=2"
